# Bandhan Mutual Fund portfolio sheet rework:
#  - lower-case / rename several headers
#  - insert a new empty "coupon" column, de-duplicate the old Yield/Quantity
#    columns down to a single "quantity" column, add three new yield-related
#    columns, and move Type/Scheme/AmcName to the end (K/L/M)
#  - refresh the "Type" column values with the new "NAN nan..." suffix
#  - extend the sheet's used range out to column M

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nRows = 4

# ---------------------------------------------------------------------
# 1. Capture the old data (rows 2-4, columns A-J) before anything moves.
# ---------------------------------------------------------------------
$old = @{}
for ($r = 2; $r -le $nRows; $r++) {
    $old[$r] = @{}
    for ($c = 1; $c -le 10; $c++) {
        $old[$r][$c] = $ws.Cells.Item($r, $c).Value()
    }
}

# New "Yield to Maturity (YTM)" values that the refreshed parser now emits
# in column J, one per data row.
$ytm = @{ 2 = "0.067293"; 3 = "0.067426"; 4 = "0.067116" }

# A lone apostrophe assigns an explicit empty TEXT value (matches the
# source's blank-but-text-typed cells) instead of leaving a truly blank
# cell behind.
$blank = "'"

# ---------------------------------------------------------------------
# 2. Clear all existing cell contents (formatting is left alone for now;
#    it gets rebuilt from scratch in step 5 below).
# ---------------------------------------------------------------------
$ws.Cells.ClearContents()

# ---------------------------------------------------------------------
# 3. Write the new header row (A1:M1), lower-cased / renamed per the
#    refreshed column schema. Every cell is written as TEXT (flip the
#    number format to "@" for the write, then straight back to General)
#    so numeric-looking values stay strings instead of becoming numbers.
# ---------------------------------------------------------------------
$headers = @{
    1  = "name of instrument"
    2  = "isin"
    3  = "coupon"
    4  = "industry"
    5  = "quantity"
    6  = "market value (mkt)"
    7  = "% to net assets (nav)"
    8  = "yield"
    9  = "yield to call (ytc)"
    10 = "yield to maturity (ytm)"
    11 = "Type"
    12 = "Scheme"
    13 = "AmcName"
}

for ($c = 1; $c -le 13; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.NumberFormat = "@"
    $cell.Value2 = $headers[$c]
    $cell.NumberFormat = "General"
}

# ---------------------------------------------------------------------
# 4. Write the data rows using the new column layout:
#      A name of instrument      <- old A
#      B isin                    <- old B
#      C coupon                  <- (new, blank)
#      D industry                <- old C
#      E quantity                <- old E (de-duplicated old D/E)
#      F market value (mkt)      <- old F
#      G % to net assets (nav)   <- old G
#      H yield                   <- (new, blank)
#      I yield to call (ytc)     <- (new, blank)
#      J yield to maturity (ytm) <- (new)
#      K Type                    <- old H, value refreshed
#      L Scheme                  <- old I
#      M AmcName                 <- old J
# ---------------------------------------------------------------------
for ($r = 2; $r -le $nRows; $r++) {
    $rowVals = $old[$r]

    $map = @{
        1  = $rowVals[1]
        2  = $rowVals[2]
        3  = $blank
        4  = $rowVals[3]
        5  = $rowVals[5]
        6  = $rowVals[6]
        7  = $rowVals[7]
        8  = $blank
        9  = $blank
        10 = $ytm[$r]
        11 = ($rowVals[8].Trim() + "  NAN nan nan nan nan nan")
        12 = $rowVals[9]
        13 = $rowVals[10]
    }

    for ($c = 1; $c -le 13; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value2 = $map[$c]
        $cell.NumberFormat = "General"
    }
}

# ---------------------------------------------------------------------
# 5. Re-apply clean formatting across the whole rebuilt range. This is
#    done LAST (copy/paste-format only touches formatting, not values)
#    so it both (a) stretches the original bold/boxed header style and
#    plain data style across the widened A1:M4 range, and (b) wipes out
#    any stray "quote prefix" / custom-number-format marks picked up
#    while coercing values to text in steps 3-4 above.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A2").Copy()
$ws.Range("A2:M4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = 0
